$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.831.76'
$ws.Range('E2').Value = '  -0.32%  '
$ws.Range('D3').Value = '2.676.30'
$ws.Range('E3').Value = '  -0.53%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = "'" + '600.32'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.97%  '
$ws.Range('D6').Value = "'" + '157.69'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.45%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  +3.51%  '
$ws.Range('E9').Value = '  +2.57%  '
$ws.Range('E10').Value = '  -0.86%  '
$ws.Range('E11').Value = '  -2.69%  '
$ws.Range('E13').Value = '  -3.78%  '
$ws.Range('D14').Value = "'" + '29.16'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.45%  '
$ws.Range('D15').Value = '3.156.11'
$ws.Range('E15').Value = '  -0.63%  '
$ws.Range('D16').Value = '65.668.24'
$ws.Range('E16').Value = '  -0.38%  '
$ws.Range('D17').Value = '2.663.36'
$ws.Range('E17').Value = '  -1.06%  '
$ws.Range('D18').Value = "'" + '12.78'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.19%  '
$ws.Range('D19').Value = "'" + '4.81'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.62%  '
$ws.Range('E20').Value = '  -3.77%  '
$ws.Range('D21').Value = "'" + '352.69'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.67%  '
$ws.Range('D22').Value = "'" + '0.999'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Value = "'" + '69.58'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.53%  '
$ws.Range('E24').Value = '  +3.34%  '
$ws.Range('D25').Value = "'" + '9.75'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.43%  '
$ws.Range('D26').Value = "'" + '1.67'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.34%  '
$ws.Range('D27').Value = "'" + '1.60'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.70%  '
$ws.Range('E28').Value = '  -3.73%  '
$ws.Range('D29').Value = "'" + '8.04'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.02%  '
$ws.Range('E30').Value = '  -0.33%  '
$ws.Range('E31').Value = '  -3.34%  '
$ws.Range('D32').Value = "'" + '535.05'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.84%  '
$ws.Range('E33').Value = '  -1.23%  '
$ws.Range('D34').Value = "'" + '6.51'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.82%  '
$ws.Range('E35').Value = '  -0.32%  '
$ws.Range('E36').Value = '  -2.98%  '
$ws.Range('E37').Value = '  -0.73%  '
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('D39').Value = "'" + '158.35'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.89%  '
$ws.Range('E40').Value = '  -2.48%  '
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('D42').Value = "'" + '163.40'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.03%  '
$ws.Range('E43').Value = '  -1.54%  '
$ws.Range('E44').Value = '  +2.45%  '
$ws.Range('E45').Value = '  -2.53%  '
$ws.Range('E46').Value = '  -4.48%  '
$ws.Range('E47').Value = '  -3.54%  '
$ws.Range('E48').Value = '  -2.64%  '
$ws.Range('D49').Value = '0.0₆0257'
$ws.Range('E49').Value = '  +8.66%  '
$ws.Range('D50').Value = "'" + '20.10'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.41%  '
$ws.Range('D51').Value = "'" + '0.0989'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.60%  '

Write-Host "Applied cryptos update."
